$d = $word.ActiveDocument

# --- Paragraph 5: "關鍵詞 (Biblica) (Chinese (Traditional)) is based on: Biblica Bible
# Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license."
# becomes the Biblica Study Notes resource blurb. Remove both hyperlinks first (deleting
# a hyperlink's own Range keeps the surrounding plain-text offsets sane; Find ranges that
# span into/out of a w:hyperlink mis-delete because of hidden field-code characters).
while ($d.Hyperlinks.Count -gt 0) {
    $h = $d.Hyperlinks.Item($d.Hyperlinks.Count)
    $hrng = $d.Range($h.Range.Start, $h.Range.End)
    $hrng.Delete()
}

$p5 = $d.Paragraphs.Item(5)
$rng = $d.Range($p5.Range.Start, $p5.Range.End)
$rng.Find.Execute("關鍵詞 (Biblica)", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Biblica Study Notes (Key Terms)", 2)

$p5 = $d.Paragraphs.Item(5)
$rng = $d.Range($p5.Range.Start, $p5.Range.End)
$rng.Find.Execute(" (Chinese (Traditional)) is based on", $false, $false, $false, $false, $false, $true, 1, $false, `
    " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. ", 2)

$p5 = $d.Paragraphs.Item(5)
$rng = $d.Range($p5.Range.Start, $p5.Range.End)
$rng.Find.Execute(": Biblica Bible Dictionary, , 2023, which is licensed under a .", $false, $false, $false, $false, $false, $true, 1, $false, `
    "Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.", 2)

# --- Remove the italic "奉耶穌的名" paragraph that sits right after the "feng" heading.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "奉耶穌的名`r") {
        $p.Range.Delete()
        break
    }
}

# --- Remove the "This PDF version is provided under the same license." paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "This PDF version is provided under the same license.`r") {
        $p.Range.Delete()
        break
    }
}

# --- Remove the "License Information" heading paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "License Information`r") {
        $p.Range.Delete()
        break
    }
}
